$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(1.182963722358093, 0.8953231288065846, 0.719691949601007, 3.4631591416293097)"
$ws.Range("C2").Value = "JSU(-1.580640448169203, 1.3100630563247049, 1.645029330324268, 4.694478264812261)"
$ws.Range("D2").Value = "JSU(-0.7804870574961984, 0.9545331454447986, 1.332053922049273, 2.0885774360867373)"
$ws.Range("E2").Value = "NIG(1.5871760264598676, 1.112123289307422, 4.437885078956225, 6.712673246740936)"
